# Adds newly crawled #belajaronline tweets (rows 198-213) to Sheet1,
# mirroring the upstream data_crawling notebook re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @()

$row198 = @{
    RowNum = 198
    A = 1346400088969543936
    B = @'
Tahun baru masih anget..5 Januari 2021. Persiapan Daring lg.
#newyear
#belajaronline
#menantivaksin
#covid19
#bekasi
'@
    C = @'
dannyurban
'@
    D = @'
Tue Jan 05 10:16:10 +0000 2021
'@
}
$newRows += $row198

$row199 = @{
    RowNum = 199
    A = 1346321621141872896
    B = @'
Last but not least, asik!
Terus perluas jaringan dan kenal sama orang baru biar tau kisah-kisah keberhasilan mereka. Et, tapi tetap harus optimis juga sama diri sendiri.
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja https://t.co/t9tDGJLoCe
'@
    C = @'
kelaskitadotcom
'@
    D = @'
Tue Jan 05 05:04:22 +0000 2021
'@
}
$newRows += $row199

$row200 = @{
    RowNum = 200
    A = 1346320666920795904
    B = @'
Emang butuh pede yang besar buat ngubah sesuatu dan belajar hal baru, tapi ya harus pede! 
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja https://t.co/55WbqzVRLF
'@
    C = @'
kelaskitadotcom
'@
    D = @'
Tue Jan 05 05:00:35 +0000 2021
'@
}
$newRows += $row200

$row201 = @{
    RowNum = 201
    A = 1346320324917296896
    B = @'
Harus berani cari tantangan baru, ya namanya juga keluar zona nyaman~
Evaluasi juga jangan ditinggalin.
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja https://t.co/8O6hMK65QF
'@
    C = @'
kelaskitadotcom
'@
    D = @'
Tue Jan 05 04:59:13 +0000 2021
'@
}
$newRows += $row201

$row202 = @{
    RowNum = 202
    A = 1346320063717016064
    B = @'
Mau ngumpulin orang-orang yang lagi bangun niat buat keluar dari zona nyaman, nih!
Simak tipsnya yuk!
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja https://t.co/FZfGa6PoVy
'@
    C = @'
kelaskitadotcom
'@
    D = @'
Tue Jan 05 04:58:11 +0000 2021
'@
}
$newRows += $row202

$row203 = @{
    RowNum = 203
    A = 1346248202039640064
    B = @'
Promo layar Interaktif Flat Panel ukuran layar lebar 65 inch ... Hub. 081386785850 @mandiriartha #korporasijabar #mandiriarthasolusi  #interaktifflatpanel #pembelajaranjarakjauh #belajaronline #bekerjajarakjauh #wfh https://t.co/rUQ11W3LWv
'@
    C = @'
echo_kilo_oscar
'@
    D = @'
Tue Jan 05 00:12:38 +0000 2021
'@
}
$newRows += $row203

$row204 = @{
    RowNum = 204
    A = 1346124985052499968
    B = @'
Pergulatan peronlinenan ummat dimulai..
kau membentukku semakin tangguh..
semangat berproses kembali dan semoga cerita yang kita gores lebih berwarna nan bermakna..
#seninvibbes #postivvibes #semangatsenin #belajaronline #dirumahsaja https://t.co/mUQ8gsnUDC
'@
    C = @'
sofi_ksmnrsh
'@
    D = @'
Mon Jan 04 16:03:00 +0000 2021
'@
}
$newRows += $row204

$row205 = @{
    RowNum = 205
    A = 1346088387938139904
    B = @'
Pemerintah pusat mengizinkan pemerintah daerah untuk melaksanakan pembelajaran tatap muka (PTM).
#BelajarDariRumah #belajaronline #daring #KBM #UPI #IKAUPI #lampung #lampostco
https://t.co/rb8U9YjcY0
'@
    C = @'
lampostco
'@
    D = @'
Mon Jan 04 13:37:35 +0000 2021
'@
}
$newRows += $row205

$row206 = @{
    RowNum = 206
    A = 1346067103321739008
    B = @'
Halo sobat Homelab
Jadikan kegagalan itu sebuah bukti kamu pernah berjuang dan jangan berhenti disituasi itu tapi selesaikan perjuanganmu sampai tuntas dan berhasil!
#homelab #askhomelab #belajarbarengberkolaborasibareng #belajardimanasaja #elearning #belajaronline https://t.co/4LT3plS06D
'@
    C = @'
Homelabmedia
'@
    D = @'
Mon Jan 04 12:13:00 +0000 2021
'@
}
$newRows += $row206

$row207 = @{
    RowNum = 207
    A = 1346042820201598976
    B = @'
Dinas Pendidikan Kabupaten Pesawaran, memutuskan kegiatan belajar mengajar (KBM) tatap muka di Bumi Andan Jejama diundur sampai waktu yang belum ditentukan.
#BelajarDariRumah #belajaronline #daring #KBM #pesawaran #lampung #lampostco
https://t.co/SQduyCyQxr
'@
    C = @'
lampostco
'@
    D = @'
Mon Jan 04 10:36:31 +0000 2021
'@
}
$newRows += $row207

$row208 = @{
    RowNum = 208
    A = 1345988381856076032
    B = @'
Makan sekuteng di rawa-rawa
Kelaskita dateng bawa kabar gembira!
Akhirnya pantun go Internasional, gais!
#kelaskita #carabarubelajarseru #belajardirumah #elearning #belajaronline #dirumahaja #taugaksih https://t.co/cJ6ld6Yok4
'@
    C = @'
kelaskitadotcom
'@
    D = @'
Mon Jan 04 07:00:12 +0000 2021
'@
}
$newRows += $row208

$row209 = @{
    RowNum = 209
    A = 1345957551767126016
    B = @'
Hi Fellas! Tak perlu repot-repot pergi ke lembaga kursus, Lister siap membantu Kamu jadi lebih hebat di tahun 2021. 
#belajaronline #lister #kursusielts #kursustoelf #kursusbahasaasing #kursusbimbinganbeasiswa #listerlanguagemaster https://t.co/QeToE9CDb4
'@
    C = @'
lister_id
'@
    D = @'
Mon Jan 04 04:57:41 +0000 2021
'@
}
$newRows += $row209

$row210 = @{
    RowNum = 210
    A = 1345908839414071040
    B = @'
Kion Kids mengandung Ion Negatif &amp; FIR yang bantu jaga kesehatan mata siBuah Hati.
#kionnanokids
#kacamatakesehatan
#belajaronline
#inspiradzi
#digitalnetworkmarketing
'@
    C = @'
SusantoSane
'@
    D = @'
Mon Jan 04 01:44:07 +0000 2021
'@
}
$newRows += $row210

$row211 = @{
    RowNum = 211
    A = 1345496765936455936
    B = @'
Jadwal Baru Belajar dari Rumah TVRI untuk PAUD dan SD Kelas 1-6, Tiap Senin-Jumat Selama 30 Menit
https://t.co/kk7OLxdr7G #JadwalBaru #BelajarDariRumah #TVRI #BelajarOnline
'@
    C = @'
tribunkaltim
'@
    D = @'
Sat Jan 02 22:26:41 +0000 2021
'@
}
$newRows += $row211

$row212 = @{
    RowNum = 212
    A = 1345339053508414976
    B = @'
Ada yang suka nonton Jdrama? Kali ini kita bakal main tebak-tebakan tentang wajah aktris atau aktor Jepang dari potongan fotonya aja nih!

#belajarbahasajepang #bahasajepang #belajaronline #nihongo #edukasi
#belajarbahasajepang #bahasajepang #belajaronline #nihongo #edukasi https://t.co/C4SLoRFGzk
'@
    C = @'
AkiNoSoraID
'@
    D = @'
Sat Jan 02 12:00:00 +0000 2021
'@
}
$newRows += $row212

$row213 = @{
    RowNum = 213
    A = 1345216933696705024
    B = @'
Yuk..yuk yang mau belajar public speaking bisa subs dan tonton videonya di link berikut https://t.co/oVBKsAw2bF 
#publicspeaking #belajaronline https://t.co/fesKWWClOg
'@
    C = @'
andigarmadi
'@
    D = @'
Sat Jan 02 03:54:44 +0000 2021
'@
}
$newRows += $row213

foreach ($r in $newRows) {
    $rowNum = $r.RowNum
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Rows.Item($rowNum).AutoFit()
}

$lastRow = $newRows[$newRows.Length - 1].RowNum

# Reflect the scroll/selection position left behind after pasting the new rows
$win = $excel.ActiveWindow
$win.ScrollRow = 193
$win.ScrollColumn = 1
[void]$ws.Range("J206").Select()

Write-Host "Added" $newRows.Length "rows through row" $lastRow
